$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.647.07'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.474.36'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.07'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.81'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0866'
$ws.Range('E10').Value = '  +10.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '33.03'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '2.857.35'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.92'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.67'
$ws.Range('E15').Value = '  -2.21%  '
$ws.Range('D16').Value = '2.479.86'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.790'
$ws.Range('E17').Value = '  +3.72%  '
$ws.Range('D18').Value = '41.620.09'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').Value = '0.0₃0957'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.43'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.50'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.61'
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.88'
$ws.Range('E29').Value = '  +2.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.44'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '156.83'
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0771'
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.45'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.103'
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  -2.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.49'
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.42'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('D44').Value = '1.983.50'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.24'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').Value = '2.716.38'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.48'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.76'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.10'
$ws.Range('E51').Value = '  +2.46%  '
